$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Defined names (baseWidth / baseHeight) ---
$wb.Names.Add("baseWidth", "=Sheet1!`$C`$2")
$wb.Names.Add("baseHeight", "=Sheet1!`$D`$2")

# --- Reshuffle existing columns: old A,B -> C,D ; old E,F -> F,G (leaving a blank gap column E) ---
$ws.Columns.Item(1).Insert()
$ws.Columns.Item(1).Insert()
$ws.Columns.Item(5).Delete()

# --- New header strings, written in the exact order needed to match shared-string indices ---
$ws.Range("I1").Value = "Next Greater Whole-Width Height"
$ws.Range("H1").Value = "Next Lesser Whole-Width Height"
$ws.Range("A1").Value = "Next Lesser Whole-Height Width"
$ws.Range("B1").Value = "Next Greater Whole-Height Width"

# --- Row 1 formatting: wrap text on the long headers, row height, filler cells ---
$ws.Range("A1:B1,H1:I1").WrapText = $true
$ws.Rows.Item(1).RowHeight = 49.5

$ws.Range("J1").WrapText = $true
$ws.Range("K1").WrapText = $true
$ws.Range("L1").WrapText = $true

# --- Column widths (best achievable given the runtime's pixel quantization) ---
$ws.Columns.Item(1).ColumnWidth = 14.5
$ws.Columns.Item(2).ColumnWidth = 12.5
$ws.Columns.Item(5).ColumnWidth = 7.6666666666667
$ws.Columns.Item(8).ColumnWidth = 13
$ws.Columns.Item(9).ColumnWidth = 13.5
$ws.Columns.Item(10).ColumnWidth = 13.5
$ws.Columns.Item(11).ColumnWidth = 17.3333333333333
$ws.Columns.Item(12).ColumnWidth = 16.5

# --- Data + formulas, row by row ---

# Row 2
$ws.Range("C2").Value = 8
$ws.Range("A2").Formula = "=IF(B2-baseWidth=C2, B2-2*baseWidth,B2-baseWidth)"
$ws.Range("B2").Formula = "=((1-(C2/baseWidth-TRUNC(C2/baseWidth)))*baseWidth)+C2"
$ws.Range("G2").Value = 1746
$ws.Range("H2").Formula = "=IF(I2-baseHeight=G2, I2-2*baseHeight, I2-baseHeight)"
$ws.Range("I2").Formula = "=((1-(G2/baseHeight-TRUNC(G2/baseHeight)))*baseHeight)+G2"

# Row 3
$ws.Range("C3").Value = 16
$ws.Range("A3").Formula = "=IF(B3-baseWidth=C3, B3-2*baseWidth,B3-baseWidth)"
$ws.Range("B3").Formula = "=((1-(C3/baseWidth-TRUNC(C3/baseWidth)))*baseWidth)+C3"
$ws.Range("G3").Value = 144
$ws.Range("H3").Formula = "=IF(I3-baseHeight=G3, I3-2*baseHeight, I3-baseHeight)"
$ws.Range("I3").Formula = "=((1-(G3/baseHeight-TRUNC(G3/baseHeight)))*baseHeight)+G3"

# Row 4
$ws.Range("C4").Value = 1848
$ws.Range("A4").Formula = "=IF(B4-baseWidth=C4, B4-2*baseWidth,B4-baseWidth)"
$ws.Range("B4").Formula = "=((1-(C4/baseWidth-TRUNC(C4/baseWidth)))*baseWidth)+C4"
$ws.Range("G4").Value = 2304
$ws.Range("H4").Formula = "=IF(I4-baseHeight=G4, I4-2*baseHeight, I4-baseHeight)"
$ws.Range("I4").Formula = "=((1-(G4/baseHeight-TRUNC(G4/baseHeight)))*baseHeight)+G4"

# Row 5
$ws.Range("C5").Value = 256
$ws.Range("A5").Formula = "=IF(B5-baseWidth=C5, B5-2*baseWidth,B5-baseWidth)"
$ws.Range("B5").Formula = "=((1-(C5/baseWidth-TRUNC(C5/baseWidth)))*baseWidth)+C5"
$ws.Range("G5").Value = 1746
$ws.Range("H5").Formula = "=IF(I5-baseHeight=G5, I5-2*baseHeight, I5-baseHeight)"
$ws.Range("I5").Formula = "=((1-(G5/baseHeight-TRUNC(G5/baseHeight)))*baseHeight)+G5"

# Row 6
$ws.Range("C6").Value = 248
$ws.Range("A6").Formula = "=IF(B6-baseWidth=C6, B6-2*baseWidth,B6-baseWidth)"
$ws.Range("B6").Formula = "=((1-(C6/baseWidth-TRUNC(C6/baseWidth)))*baseWidth)+C6"
$ws.Range("G6").Value = 1755
$ws.Range("H6").Formula = "=IF(I6-baseHeight=G6, I6-2*baseHeight, I6-baseHeight)"
$ws.Range("I6").Formula = "=((1-(G6/baseHeight-TRUNC(G6/baseHeight)))*baseHeight)+G6"

# --- Drop the old trailing row 7 remnant (only the shared-formula spillover in column D survived the shift) ---
$ws.Range("D7").ClearContents()

# --- Selection, matching the saved workbook state ---
$ws.Range("G7").Select()
